$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.215.23'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '2.384.28'
$ws.Range('E3').Value = '  -3.66%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '549.07'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = '141.22'
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -10.63%  '
$ws.Range('D9').Value = '2.383.75'
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('D14').Value = '25.45'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').Value = '2.817.13'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '60.774.42'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '2.386.53'
$ws.Range('E18').Value = '  -3.63%  '
$ws.Range('D19').Value = '10.75'
$ws.Range('E19').Value = '  -3.77%  '
$ws.Range('D20').Value = '4.12'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').Value = '318.18'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').Value = '6.70'
$ws.Range('E22').Value = '  -4.95%  '
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '63.43'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('D26').Value = '8.25'
$ws.Range('E26').Value = '  +5.77%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '2.504.60'
$ws.Range('E28').Value = '  -3.93%  '
$ws.Range('D29').Value = '0.0₃0927'
$ws.Range('E29').Value = '  -6.44%  '
$ws.Range('D30').Value = '523.15'
$ws.Range('E30').Value = '  -3.12%  '
$ws.Range('E31').Value = '  -5.06%  '
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('E33').Value = '  -3.98%  '
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').Value = '5.51'
$ws.Range('E37').Value = '  -6.19%  '
$ws.Range('D38').Value = '4.68'
$ws.Range('E38').Value = '  -4.00%  '
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E40').Value = '  +6.48%  '
$ws.Range('D41').Value = '18.07'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').Value = '139.31'
$ws.Range('E42').Value = '  -4.65%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = '2.15'
$ws.Range('E45').Value = '  -9.14%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '3.63'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '140.52'
$ws.Range('E47').Value = '  -4.68%  '
$ws.Range('D48').Value = '20.15'
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('E49').Value = '  -3.75%  '
$ws.Range('E50').Value = '  -3.68%  '
$ws.Range('E51').Value = '  -0.85%  '
